# Auto-generated edit script for uniao_da_vitoria.xlsx update (27-10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-63 with new computed statistics ---
$ws.Range("F2").Value = 0.456438285308356
$ws.Range("G2").Value = 0.5946080815016426
$ws.Range("F3").Value = 0.1188904404129971
$ws.Range("G3").Value = 0.292207416994111
$ws.Range("F4").Value = 0.04934630657270472
$ws.Range("G4").Value = 0.1859492534638571
$ws.Range("F5").Value = 0.03996727157481791
$ws.Range("G5").Value = 0.1669321499764016
$ws.Range("F6").Value = 0.07849757518135056
$ws.Range("G6").Value = 0.2359414138669834
$ws.Range("F7").Value = 0.4546940260524142
$ws.Range("G7").Value = 0.5933556422553211
$ws.Range("F8").Value = 0.001016679896926618
$ws.Range("G8").Value = 0.02568731204238746
$ws.Range("F9").Value = 0.5331444704134104
$ws.Range("G9").Value = 0.6483240132782805
$ws.Range("F10").Value = 0.3810327115690741
$ws.Range("G10").Value = 0.5388323385177153
$ws.Range("F11").Value = 0.108607523417468
$ws.Range("G11").Value = 0.2788730313709404
$ws.Range("D12").Value = 41.64799999999998
$ws.Range("F12").Value = 0.05415296938776069
$ws.Range("G12").Value = 0.1950169135912066
$ws.Range("F13").Value = 0.03028856132912447
$ws.Range("G13").Value = 0.1448698362448913
$ws.Range("F14").Value = 0.009317928295665977
$ws.Range("G14").Value = 0.07939398569248958
$ws.Range("F15").Value = 0.2029488327451937
$ws.Range("G15").Value = 0.3857813671399947
$ws.Range("F16").Value = 0.948769130448873
$ws.Range("G16").Value = 0.9359366468609942
$ws.Range("D17").Value = 697.0880999999997
$ws.Range("F17").Value = 0.8073484016113568
$ws.Range("G17").Value = 0.8309037999330935
$ws.Range("F18").Value = 0.8628971389059762
$ws.Range("G18").Value = 0.869406435400824
$ws.Range("A19").Value = 38
$ws.Range("F19").Value = 0.3635048544208712
$ws.Range("G19").Value = 0.5253095515545443
$ws.Range("A20").Value = 40
$ws.Range("F20").Value = 0.8514897811709355
$ws.Range("G20").Value = 0.8613344367964542
$ws.Range("A21").Value = 44
$ws.Range("D21").Value = 370.7196999999999
$ws.Range("F21").Value = 0.5745273194679344
$ws.Range("G21").Value = 0.6764188183050507
$ws.Range("A22").Value = 49
$ws.Range("D22").Value = 1352.8269
$ws.Range("F22").Value = 0.9607888661344086
$ws.Range("G22").Value = 0.9468712271283568
$ws.Range("A23").Value = 52
$ws.Range("F23").Value = 0.1948010493038771
$ws.Range("G23").Value = 0.3776063161124989
$ws.Range("A24").Value = 54
$ws.Range("F24").Value = 0.8621092778049834
$ws.Range("G24").Value = 0.8688454646257577
$ws.Range("A25").Value = 59
$ws.Range("F25").Value = 0.07594396051253505
$ws.Range("G25").Value = 0.23196556401993
$ws.Range("A26").Value = 62
$ws.Range("F26").Value = 0.495298089818908
$ws.Range("G26").Value = 0.6221299527131881
$ws.Range("A27").Value = 68
$ws.Range("D27").Value = 938.4727999999998
$ws.Range("F27").Value = 0.8927804017850125
$ws.Range("G27").Value = 0.8911477512829834
$ws.Range("A28").Value = 70
$ws.Range("F28").Value = 0.1545605211018792
$ws.Range("G28").Value = 0.3347399160960936
$ws.Range("A29").Value = 72
$ws.Range("F29").Value = 0.2244693258856998
$ws.Range("G29").Value = 0.4067000445763652
$ws.Range("A30").Value = 76
$ws.Range("F30").Value = 0.9834109810777103
$ws.Range("G30").Value = 0.9705449141039382
$ws.Range("A31").Value = 77
$ws.Range("F31").Value = 0.3991451284445828
$ws.Range("G31").Value = 0.5525639154786957
$ws.Range("A32").Value = 79
$ws.Range("F32").Value = 0.6484864247347032
$ws.Range("G32").Value = 0.7256515131078404
$ws.Range("A33").Value = 80
$ws.Range("F33").Value = 0.03267108690801651
$ws.Range("G33").Value = 0.1505850884739003
$ws.Range("A34").Value = 81
$ws.Range("F34").Value = 0.1317573562949799
$ws.Range("G34").Value = 0.3081553789870479
$ws.Range("A35").Value = 84
$ws.Range("C35").Value = 25137
$ws.Range("D35").Value = 870.3554999999998
$ws.Range("E35").Value = 88
$ws.Range("F35").Value = 0.8734995432028833
$ws.Range("G35").Value = 0.8770109888807145
$ws.Range("A36").Value = 86
$ws.Range("F36").Value = 0.4708679492544749
$ws.Range("G36").Value = 0.6049103722871053
$ws.Range("A37").Value = 89
$ws.Range("F37").Value = 0.365957456839671
$ws.Range("G37").Value = 0.5272164035966458
$ws.Range("A38").Value = 90
$ws.Range("D38").Value = 193.402
$ws.Range("F38").Value = 0.3456312613568117
$ws.Range("G38").Value = 0.5112601163553028
$ws.Range("A39").Value = 93
$ws.Range("F39").Value = 0.366785033373981
$ws.Range("G39").Value = 0.5278587260426937
$ws.Range("A40").Value = 95
$ws.Range("F40").Value = 0.6197667320535365
$ws.Range("G40").Value = 0.706644830267916
$ws.Range("A41").Value = 96
$ws.Range("F41").Value = 0.007008550576576805
$ws.Range("G41").Value = 0.06866786639844855
$ws.Range("A42").Value = 100
$ws.Range("F42").Value = 0.9794497154816091
$ws.Range("G42").Value = 0.9659376638389008
$ws.Range("A43").Value = 104
$ws.Range("F43").Value = 0.2102921147728999
$ws.Range("G43").Value = 0.3930248707346073
$ws.Range("A44").Value = 106
$ws.Range("F44").Value = 0.1212628699479383
$ws.Range("G44").Value = 0.295206266264294
$ws.Range("A45").Value = 109
$ws.Range("F45").Value = 0.04279019191667693
$ws.Range("G45").Value = 0.1728636118895857
$ws.Range("A46").Value = 110
$ws.Range("F46").Value = 0.2703000266792889
$ws.Range("G46").Value = 0.4485209877096592
$ws.Range("A47").Value = 115
$ws.Range("D47").Value = 3052.843299999999
$ws.Range("F47").Value = 0.9993675357725381
$ws.Range("G47").Value = 0.9965394875327398
$ws.Range("A48").Value = 117
$ws.Range("F48").Value = 0.0825772855254862
$ws.Range("G48").Value = 0.2421677103570429
$ws.Range("A49").Value = 119
$ws.Range("F49").Value = 0.02753501549482289
$ws.Range("G49").Value = 0.1379852508818057
$ws.Range("A50").Value = 120
$ws.Range("F50").Value = 0.2718528427878296
$ws.Range("G50").Value = 0.4498823105983869
$ws.Range("A51").Value = 121
$ws.Range("F51").Value = 0.1171745368948168
$ws.Range("G51").Value = 0.2900209881439849
$ws.Range("A52").Value = 122
$ws.Range("F52").Value = 0.07771514914858633
$ws.Range("G52").Value = 0.2347298494939028
$ws.Range("A53").Value = 125
$ws.Range("F53").Value = 0.5995345463383216
$ws.Range("G53").Value = 0.6931786664959569
$ws.Range("A54").Value = 127
$ws.Range("F54").Value = 0.4367896741435222
$ws.Range("G54").Value = 0.5804058927496979
$ws.Range("A55").Value = 128
$ws.Range("F55").Value = 0.02823357069298679
$ws.Range("G55").Value = 0.1397623444467889
$ws.Range("A56").Value = 132
$ws.Range("B56").Value = 34620
$ws.Range("D56").Value = 139.9039999999999
$ws.Range("E56").Value = 7
$ws.Range("F56").Value = 0.2548787708821944
$ws.Range("G56").Value = 0.4348164264292126
$ws.Range("A57").Value = 133
$ws.Range("G57").Value = [double]"1.759430566470838e-09"
$ws.Range("A58").Value = 137
$ws.Range("F58").Value = 0.3315888457356055
$ws.Range("G58").Value = 0.5000217193025353
$ws.Range("A59").Value = 140
$ws.Range("F59").Value = 0.5756241970409621
$ws.Range("G59").Value = 0.6771570090182866
$ws.Range("A60").Value = 145
$ws.Range("F60").Value = 0.9862774060705517
$ws.Range("G60").Value = 0.9740845997455293
$ws.Range("A61").Value = 154
$ws.Range("F61").Value = 0.671434627507651
$ws.Range("G61").Value = 0.7407765438013726
$ws.Range("F62").Value = 0.09194020613779094
$ws.Range("G62").Value = 0.255926962991087
$ws.Range("F63").Value = 0.9997699366975298
$ws.Range("G63").Value = 0.9981836187804234

# --- Append new row 64 (period id 164) ---
# Copy formatting (styles) from row 63 so the new row matches existing look
$ws.Range("A63:G63").Copy() | Out-Null
$ws.Range("A64:G64").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A64").Value = 164
$ws.Range("B64").Value = 44123
$ws.Range("C64").Value = 44130
$ws.Range("D64").Value = 318.822396
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = 0.5173987797657653
$ws.Range("G64").Value = 0.6374912600607676
